$d = $word.ActiveDocument

# Target: the signature table (7th table in the document) contains a row
# with an empty cell that has 4 empty "Encabezado"-styled paragraphs.
# The edit removes the first 3 of those empty paragraphs, leaving just 1.
$t = $d.Tables.Item(7)
$cell = $t.Cell(2, 2)

# Delete the first three paragraphs in the cell, leaving only the last one.
for ($i = 1; $i -le 3; $i++) {
    $cell.Range.Paragraphs.Item(1).Range.Delete()
}
